$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header title (A1) to reflect new milestone wording
$ws.Range("A1").Value = "contributions for we move"

# Add the new contributor entry in row 4 (A4), matching the trailing space in the source text
$ws.Range("A4").Value = "kian "

# Move the active selection to B4 to match the saved view state
$ws.Range("B4").Select()
